$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update adds a new week's pair of rows ("Primera" / "Segunda"
# quality grades) at the top of this market/category data block. The
# existing rows 263..276 all shift down by two (to 265..278), and the two
# new rows are inserted at 263/264.

$ws.Range("A263:R264").EntireRow.Insert()

# Row 263 - Zanahoria, Primera
$ws.Cells.Item(263,1).Value  = 11
$ws.Cells.Item(263,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(263,3).Value  = "Bíobío"
$ws.Cells.Item(263,4).Value  = 44826
$ws.Cells.Item(263,5).Value  = 8
$ws.Cells.Item(263,6).Value  = 100114013
$ws.Cells.Item(263,7).Value  = "Zanahoria"
$ws.Cells.Item(263,8).Value  = "Sin especificar"
$ws.Cells.Item(263,9).Value  = "Primera"
$ws.Cells.Item(263,10).Value = 600
$ws.Cells.Item(263,11).Value = 9000
$ws.Cells.Item(263,12).Value = 10000
$ws.Cells.Item(263,13).Value = 9500
$ws.Cells.Item(263,14).Value = "$/saco 20 kilos"
$ws.Cells.Item(263,15).Value = "Región de La Araucanía"
$ws.Cells.Item(263,16).Value = 475
$ws.Cells.Item(263,17).Value = 20
$ws.Cells.Item(263,18).Value = "Hortaliza"

# Row 264 - Zanahoria, Segunda
$ws.Cells.Item(264,1).Value  = 11
$ws.Cells.Item(264,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(264,3).Value  = "Bíobío"
$ws.Cells.Item(264,4).Value  = 44826
$ws.Cells.Item(264,5).Value  = 8
$ws.Cells.Item(264,6).Value  = 100114013
$ws.Cells.Item(264,7).Value  = "Zanahoria"
$ws.Cells.Item(264,8).Value  = "Sin especificar"
$ws.Cells.Item(264,9).Value  = "Segunda"
$ws.Cells.Item(264,10).Value = 300
$ws.Cells.Item(264,11).Value = 8000
$ws.Cells.Item(264,12).Value = 8000
$ws.Cells.Item(264,13).Value = 8000
$ws.Cells.Item(264,14).Value = "$/saco 20 kilos"
$ws.Cells.Item(264,15).Value = "Región de La Araucanía"
$ws.Cells.Item(264,16).Value = 400
$ws.Cells.Item(264,17).Value = 20
$ws.Cells.Item(264,18).Value = "Hortaliza"
